# Add a new "exibir_como_botao" column (I) to the data sheet, mirroring
# the "exibir_ao_iniciar" column (H): header in row 1, and "sim" in the
# same rows that currently have "sim" in column H for the "fechada"/
# "fechados" (closed) state rows used as the default-shown state, plus the
# first "puxadores1_fechados" row of each hardware group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new cells the same cell style ("s=1") already used throughout
# the sheet (e.g. column A), by copying formats from A1 onto each new
# cell before filling in its value.
$ws.Range("A1").Copy()
$newCellRows = @(1, 2, 7, 12, 16, 20, 22)
foreach ($r in $newCellRows) {
    $ws.Cells.Item($r, 9).PasteSpecial(-4122)
}

# Header
$ws.Cells.Item(1, 9).Value = "exibir_como_botao"

# Data rows that get "sim" in the new column
$rows = @(2, 7, 12, 16, 20, 22)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).Value = "sim"
}

# Set the width of the new column (I) to (as close as possible to) 17.38
$ws.Columns.Item(9).ColumnWidth = 16.5
